$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03405623995107
$ws.Range("D2").Value = 1.053937298376954
$ws.Range("E2").Value = 1.033273968449183
$ws.Range("F2").Value = 1.058389320342116
$ws.Range("I2").Value = 1.041446407748981
$ws.Range("J2").Value = 1.039177576438479
$ws.Range("K2").Value = 1.05668188945192
$ws.Range("L2").Value = 1.036076790351392
$ws.Range("M2").Value = 1.06112168783682
$ws.Range("N2").Value = 1.016962485373594

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03511211875877
$ws.Range("D3").Value = 1.054711345835702
$ws.Range("E3").Value = 1.034174896856251
$ws.Range("F3").Value = 1.059367726844822
$ws.Range("I3").Value = 1.041713268513735
$ws.Range("J3").Value = 1.039875888773329
$ws.Range("K3").Value = 1.057269180782924
$ws.Range("L3").Value = 1.036786419238477
$ws.Range("M3").Value = 1.061913700999173
$ws.Range("N3").Value = 1.017199557136452

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.035795591188602
$ws.Range("D4").Value = 1.055212356676277
$ws.Range("E4").Value = 1.034758410727506
$ws.Range("F4").Value = 1.060001425775855
$ws.Range("I4").Value = 1.041884864582707
$ws.Range("J4").Value = 1.040327427848801
$ws.Range("K4").Value = 1.057648686896345
$ws.Range("L4").Value = 1.037245519337699
$ws.Range("M4").Value = 1.062426171140407
$ws.Range("N4").Value = 1.017352716064742

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.036082982452811
$ws.Range("D5").Value = 1.055423016611189
$ws.Range("E5").Value = 1.035003851726083
$ws.Range("F5").Value = 1.06026797664285
$ws.Range("I5").Value = 1.041956744525575
$ws.Range("J5").Value = 1.040517178951307
$ws.Range("K5").Value = 1.057808108297795
$ws.Range("L5").Value = 1.037438506096037
$ws.Range("M5").Value = 1.062641609157774
$ws.Range("N5").Value = 1.017417045835507

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.036131240198004
$ws.Range("D6").Value = 1.055458389362577
$ws.Range("E6").Value = 1.035045070058935
$ws.Range("F6").Value = 1.060312740115338
$ws.Range("I6").Value = 1.041968798284366
$ws.Range("J6").Value = 1.040549034534105
$ws.Range("K6").Value = 1.057834868620789
$ws.Range("L6").Value = 1.037470908299772
$ws.Range("M6").Value = 1.062677781867066
$ws.Range("N6").Value = 1.017427843665242

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.035799431091074
$ws.Range("D7").Value = 1.05521517138848
$ws.Range("E7").Value = 1.034761689805146
$ws.Range("F7").Value = 1.060004986877506
$ws.Range("I7").Value = 1.041885826064291
$ws.Range("J7").Value = 1.040329963611079
$ws.Range("K7").Value = 1.057650817576993
$ws.Range("L7").Value = 1.037248098111101
$ws.Range("M7").Value = 1.062429049852092
$ws.Range("N7").Value = 1.017353575871419

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.034413027753361
$ws.Range("D8").Value = 1.054198859588805
$ws.Range("E8").Value = 1.033578326563954
$ws.Range("F8").Value = 1.058719851704352
$ws.Range("I8").Value = 1.041536818162048
$ws.Range("J8").Value = 1.039413639521087
$ws.Range("K8").Value = 1.056880472608188
$ws.Range("L8").Value = 1.036316628490115
$ws.Range("M8").Value = 1.061389355112086
$ws.Range("N8").Value = 1.017042654933287

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031971907633337
$ws.Range("D9").Value = 1.052409185030114
$ws.Range("E9").Value = 1.031497350955861
$ws.Range("F9").Value = 1.056459947835809
$ws.Range("I9").Value = 1.040913559143294
$ws.Range("J9").Value = 1.037796556547745
$ws.Range("K9").Value = 1.055519146960509
$ws.Range("L9").Value = 1.034674687255071
$ws.Range("M9").Value = 1.059557196015681
$ws.Range("N9").Value = 1.016492926095446

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030345766534225
$ws.Range("D10").Value = 1.051216932667755
$ws.Range("E10").Value = 1.03011293331966
$ws.Range("F10").Value = 1.054956534652754
$ws.Range("I10").Value = 1.040492516389032
$ws.Range("J10").Value = 1.036716904171296
$ws.Range("K10").Value = 1.054609029641201
$ws.Range("L10").Value = 1.033579698400884
$ws.Range("M10").Value = 1.058335743458551
$ws.Range("N10").Value = 1.016125211302096

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.029641929260788
$ws.Range("D11").Value = 1.050700891300611
$ws.Range("E11").Value = 1.029514159520147
$ws.Range("F11").Value = 1.05430630673954
$ws.Range("I11").Value = 1.04030889061943
$ws.Range("J11").Value = 1.036249027276852
$ws.Range("K11").Value = 1.05421433830249
$ws.Range("L11").Value = 1.033105474881998
$ws.Range("M11").Value = 1.057806846547167
$ws.Range("N11").Value = 1.015965697704369

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029380536432043
$ws.Range("D12").Value = 1.050509243306792
$ws.Range("E12").Value = 1.029291852125524
$ws.Range("F12").Value = 1.054064898060966
$ws.Range("I12").Value = 1.040240487043985
$ws.Range("J12").Value = 1.036075180020269
$ws.Range("K12").Value = 1.054067642158573
$ws.Range("L12").Value = 1.032929314503283
$ws.Range("M12").Value = 1.057610391334591
$ws.Range("N12").Value = 1.015906403753116

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029436604109276
$ws.Range("D13").Value = 1.050550350955249
$ws.Range("E13").Value = 1.029339533109185
$ws.Range("F13").Value = 1.054116675829569
$ws.Range("I13").Value = 1.040255168744249
$ws.Range("J13").Value = 1.03611247342081
$ws.Range("K13").Value = 1.054099113051558
$ws.Range("L13").Value = 1.03296710207024
$ws.Range("M13").Value = 1.057652531605
$ws.Range("N13").Value = 1.015919124472609

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029620321558706
$ws.Range("D14").Value = 1.050685048944257
$ws.Range("E14").Value = 1.029495781400257
$ws.Range("F14").Value = 1.054286349474184
$ws.Range("I14").Value = 1.040303240371542
$ws.Range("J14").Value = 1.036234658168967
$ws.Range("K14").Value = 1.054202214190096
$ws.Range("L14").Value = 1.033090913666495
$ws.Range("M14").Value = 1.057790607481093
$ws.Range("N14").Value = 1.015960797334236

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029733521721556
$ws.Range("D15").Value = 1.050768045187944
$ws.Range("E15").Value = 1.029592064904939
$ws.Range("F15").Value = 1.054390906257046
$ws.Range("I15").Value = 1.040332832818674
$ws.Range("J15").Value = 1.036309932684731
$ws.Range("K15").Value = 1.054265726269481
$ws.Range("L15").Value = 1.033167196409819
$ws.Range("M15").Value = 1.057875680690467
$ws.Range("N15").Value = 1.015986467600763

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030392484035785
$ws.Range("D16").Value = 1.051251185174982
$ws.Range("E16").Value = 1.030152686503374
$ws.Range("F16").Value = 1.054999704212511
$ws.Range("I16").Value = 1.04050467541622
$ws.Range("J16").Value = 1.036747947623624
$ws.Range("K16").Value = 1.054635211321308
$ws.Range("L16").Value = 1.033611169264252
$ws.Range("M16").Value = 1.058370844660961
$ws.Range("N16").Value = 1.016135791584794

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030805911811829
$ws.Range("D17").Value = 1.051554303479122
$ws.Range("E17").Value = 1.030504534191916
$ws.Range("F17").Value = 1.055381791066707
$ws.Range("I17").Value = 1.040612116790759
$ws.Range("J17").Value = 1.037022601140543
$ws.Range("K17").Value = 1.054866818117417
$ws.Range("L17").Value = 1.033889638853648
$ws.Range("M17").Value = 1.058681448412934
$ws.Range("N17").Value = 1.016229380875844

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.031047085484624
$ws.Range("D18").Value = 1.051731127628309
$ws.Range("E18").Value = 1.030709827363334
$ws.Range("F18").Value = 1.055604729134824
$ws.Range("I18").Value = 1.040674658971482
$ws.Range("J18").Value = 1.037182765184883
$ws.Range("K18").Value = 1.055001852072955
$ws.Range("L18").Value = 1.034052057125074
$ws.Range("M18").Value = 1.058862618294166
$ws.Range("N18").Value = 1.016283941869457

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.0311293242741
$ws.Range("D19").Value = 1.051791423520558
$ws.Range("E19").Value = 1.030779838272485
$ws.Range("F19").Value = 1.055680757659925
$ws.Range("I19").Value = 1.040695962765445
$ws.Range("J19").Value = 1.037237370753428
$ws.Range("K19").Value = 1.055047885242844
$ws.Range("L19").Value = 1.034107436136883
$ws.Range("M19").Value = 1.058924392517083
$ws.Range("N19").Value = 1.016302540994247

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030761552004422
$ws.Range("D20").Value = 1.051521779639966
$ws.Range("E20").Value = 1.030466777396659
$ws.Range("F20").Value = 1.05534078917909
$ws.Range("I20").Value = 1.040600602436486
$ws.Range("J20").Value = 1.03699313721264
$ws.Range("K20").Value = 1.054841974933096
$ws.Range("L20").Value = 1.033859762561438
$ws.Range("M20").Value = 1.058648123568578
$ws.Range("N20").Value = 1.01621934253139

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02956622014552
$ws.Range("D21").Value = 1.050645382832319
$ws.Range("E21").Value = 1.029449767302281
$ws.Range("F21").Value = 1.054236381636462
$ws.Range("I21").Value = 1.040289089905801
$ws.Range("J21").Value = 1.036198679376466
$ws.Range("K21").Value = 1.054171855959721
$ws.Range("L21").Value = 1.033054454584163
$ws.Range("M21").Value = 1.057749947554012
$ws.Range("N21").Value = 1.015948526914563

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028814920106885
$ws.Range("D22").Value = 1.050094547297467
$ws.Range("E22").Value = 1.028810933794374
$ws.Range("F22").Value = 1.053542662043906
$ws.Range("I22").Value = 1.040092091178684
$ws.Range("J22").Value = 1.035698842827351
$ws.Range("K22").Value = 1.05375000396003
$ws.Range("L22").Value = 1.032548052457542
$ws.Range("M22").Value = 1.057185232717284
$ws.Range("N22").Value = 1.015778002801944

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029213174626349
$ws.Range("D23").Value = 1.050386537323817
$ws.Range("E23").Value = 1.029149534390942
$ws.Range("F23").Value = 1.053910352615144
$ws.Range("I23").Value = 1.040196631759034
$ws.Range("J23").Value = 1.035963846796709
$ws.Range("K23").Value = 1.053973684922509
$ws.Range("L23").Value = 1.032816512632506
$ws.Range("M23").Value = 1.057484598139138
$ws.Range("N23").Value = 1.01586842465094

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030781596204486
$ws.Range("D24").Value = 1.051536475700023
$ws.Range("E24").Value = 1.030483837862118
$ws.Range("F24").Value = 1.055359315940815
$ws.Range("I24").Value = 1.040605805668574
$ws.Range("J24").Value = 1.037006450806741
$ws.Range("K24").Value = 1.054853200679097
$ws.Range("L24").Value = 1.033873262397562
$ws.Range("M24").Value = 1.058663181631392
$ws.Range("N24").Value = 1.016223878513835

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032602771098111
$ws.Range("D25").Value = 1.052871710861837
$ws.Range("E25").Value = 1.032034824489653
$ws.Range("F25").Value = 1.057043628707788
$ws.Range("I25").Value = 1.041075663788699
$ws.Range("J25").Value = 1.038214894109675
$ws.Range("K25").Value = 1.055871537735118
$ws.Range("L25").Value = 1.035099233855066
$ws.Range("M25").Value = 1.060030858559221
$ws.Range("N25").Value = 1.016635261695945
